{"js": "// Replace the 25 two-digit-divided-by-one-digit problems in the worksheet\n// table with their new values, in document order. The table has 20 rows\n// (5 \"problem\" rows of 5 cells each, interleaved with 3 blank rows), so the\n// non-blank rows are at zero-based indices 0, 4, 8, 12, 16.\nconst oldToNew = [\n  \"19\u00f72=\", \"99\u00f77=\",\n  \"14\u00f73=\", \"17\u00f77=\",\n  \"92\u00f72=\", \"67\u00f79=\",\n  \"86\u00f75=\", \"71\u00f78=\",\n  \"48\u00f79=\", \"32\u00f73=\",\n  \"77\u00f73=\", \"10\u00f75=\",\n  \"61\u00f74=\", \"69\u00f79=\",\n  \"18\u00f75=\", \"54\u00f75=\",\n  \"27\u00f74=\", \"99\u00f73=\",\n  \"98\u00f75=\", \"95\u00f79=\",\n  \"87\u00f73=\", \"26\u00f74=\",\n  \"16\u00f74=\", \"58\u00f78=\",\n  \"47\u00f76=\", \"63\u00f75=\",\n  \"22\u00f74=\", \"81\u00f75=\",\n  \"52\u00f76=\", \"57\u00f79=\",\n  \"89\u00f72=\", \"87\u00f75=\",\n  \"17\u00f74=\", \"64\u00f74=\",\n  \"82\u00f75=\", \"19\u00f79=\",\n  \"55\u00f73=\", \"25\u00f72=\",\n  \"25\u00f73=\", \"19\u00f79=\",\n  \"39\u00f76=\", \"76\u00f74=\",\n  \"27\u00f74=\", \"94\u00f77=\",\n  \"74\u00f73=\", \"22\u00f74=\",\n  \"74\u00f79=\", \"18\u00f75=\",\n  \"38\u00f76=\", \"13\u00f72=\",\n];\nconst replacements = [];\nfor (let i = 0; i < oldToNew.length; i += 2) {\n  replacements.push(oldToNew[i + 1]);\n}\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst rows = table.values;\nlet k = 0;\nfor (let r = 0; r < rows.length && k < replacements.length; r++) {\n  const row = rows[r];\n  for (let c = 0; c < row.length && k < replacements.length; c++) {\n    const text = row[c];\n    if (text && text.indexOf(\"\u00f7\") !== -1) {\n      table.getCell(r, c).value = replacements[k];\n      k++;\n    }\n  }\n}\nawait context.sync();\n", "ps1": "# Replace the 25 two-digit-divided-by-one-digit problems in the worksheet\n# table with their new values, in document order. The table has 20 rows\n# (5 \"problem\" rows of 5 cells each, interleaved with 3 blank rows), so the\n# non-blank rows are at 1-based indices 1, 5, 9, 13, 17 \u2014 but it's simplest\n# (and robust to layout assumptions) to just walk every cell in the table in\n# reading order and replace the ones that contain a division problem.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  \"99\u00f77=\", \"17\u00f77=\", \"67\u00f79=\", \"71\u00f78=\", \"32\u00f73=\",\n  \"10\u00f75=\", \"69\u00f79=\", \"54\u00f75=\", \"99\u00f73=\", \"95\u00f79=\",\n  \"26\u00f74=\", \"58\u00f78=\", \"63\u00f75=\", \"81\u00f75=\", \"57\u00f79=\",\n  \"87\u00f75=\", \"64\u00f74=\", \"19\u00f79=\", \"25\u00f72=\", \"19\u00f79=\",\n  \"76\u00f74=\", \"94\u00f77=\", \"22\u00f74=\", \"18\u00f75=\", \"13\u00f72=\"\n)\n\n$t = $d.Tables.Item(1)\n$rows = $t.Rows.Count\n$cols = $t.Columns.Count\n\n$k = 0\nfor ($r = 1; $r -le $rows; $r++) {\n  for ($c = 1; $c -le $cols; $c++) {\n    $cell = $t.Cell($r, $c)\n    $text = $cell.Range.Text\n    if ($text -and $text.Contains(\"\u00f7\") -and $k -lt $replacements.Length) {\n      $cell.Range.Text = $replacements[$k]\n      $k = $k + 1\n    }\n  }\n}\n"}
